# CCGX-Modbus-TCP-register-list.xlsx update
# "Map unitid 228 to 275, the nanopi VE.Bus port.
#  This is for EasySolar-II-GX and Multiplus-II-GX models."

$wb = $excel.ActiveWorkbook

# --- 1. "Unit ID mapping" sheet: insert new mapping row for unit-id 228 ---
$wsMap = $wb.Worksheets.Item("Unit ID mapping")

# Insert a fresh row above the existing row 7 (pushes everything below it,
# including the merged note rows, down by one and keeps their formatting).
$wsMap.Rows.Item(7).Insert()

$wsMap.Range("A7").Value = 228
$wsMap.Range("B7").Value = 275
$wsMap.Range("C7").Value = "EasySolar-II/Multiplus-II GX VE.Bus port (ttyS3)"
$wsMap.Rows.Item(7).RowHeight = 13.8

# --- 2. "Document versions" sheet: log the change as Rev 24 ---
$wsVer = $wb.Worksheets.Item("Document versions")
$wsVer.Range("A65").Value = "Rev 24"
$wsVer.Range("B65").Value = "Added mapping for EasySolar-II/Multiplus-II GX VE.Bus port"

# --- 3. Restore view/selection state on each sheet ---
[void]$wsMap.Range("C49").Select()

[void]$wsVer.Range("B65").Select()

# "Field list" becomes the active sheet/tab again, selected last so it
# ends up as the active window/tab.
$wsField = $wb.Worksheets.Item("Field list")
[void]$wsField.Range("A3").Select()
